# Update the pl_mw.xlsx "res_line" results sheet for Case_2_246 (380 kV case).
# Only columns B, C, E, F, G, L, M, O change for data rows 2-25 (A, D, H, I, J, K, N stay the same).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that receive new values, in order.
$cols = @('B','C','E','F','G','L','M','O')

# New values for each data row (rows 2..25), one array per row, values in $cols order.
$data = @(
    @(1.139697690234698,0.2941092967180055,0.1199065955010646,0.4443680307746263,0.002446004638016376,0.1894141985854674,0.240628942188124,2.773599229613637),
    @(1.033918060399799,0.2808339296921645,0.1213875759088742,0.387822817061874,0.002448892298213491,0.1867897270564995,0.2243639549636498,2.815708847190933),
    @(0.9690008673242119,0.2726673365299632,0.122350048167013,0.3531389305168915,0.002450759036712857,0.1852698444733676,0.2144243922471105,2.843955057983592),
    @(0.9425561265253464,0.2693357656713999,0.1227556414531241,0.3390132514313251,0.002451543380652392,0.1846735488726594,0.210386037112869,2.856065624789167),
    @(0.9381656187150043,0.2687823503883351,0.1228237983344491,0.336668177824194,0.002451675050110112,0.1845759291402587,0.2097162087116757,2.858112779961573),
    @(0.9686441836462905,0.2726224200874299,0.1223554639535003,0.3529483938344953,0.002450769518830898,0.1852617091543749,0.2143698802926082,2.844115957470251),
    @(1.103219117882873,0.2895353171104205,0.1204062162634881,0.4248636149813336,0.002446980897897231,0.1884903056747547,0.2350111075121291,2.787621731956264),
    @(1.367320714952939,0.3225685823209403,0.1170047342062033,0.5661985755041457,0.002440291673589723,0.195546503624783,0.2758552700352936,2.695858156219074),
    @(1.561427767527846,0.3467459762950966,0.1147613033428757,0.6702781546542269,0.002435823820806387,0.2011714308006418,0.3060798860510303,2.640106786697203),
    @(1.649738391250366,0.3577227327257049,0.1137960023706711,0.7176906081379002,0.002433887309635268,0.2038258628883938,0.3198755503224646,2.617292814377549),
    @(1.683179620907993,0.3618760010726874,0.113438397613189,0.7356546913071611,0.002433167725483757,0.2048447449435713,0.3251060966596171,2.609021378242232),
    @(1.675977482533199,0.3609816756036253,0.113515061508435,0.7317853510981394,0.002433322091218989,0.2046247015790073,0.3239793227660925,2.61078640365821),
    @(1.652489632047889,0.358064494144628,0.113766423068256,0.7191683204515869,0.002433827834014958,0.2039094124885281,0.3203057430360303,2.616604940812266),
    @(1.638102587633796,0.3562771865125001,0.1139214221773918,0.7114413442032514,0.002434139402686206,0.2034730609094737,0.3180563995131678,2.620216891660533),
    @(1.555656533814499,0.3460281598730717,0.114825498827037,0.6671810134426437,0.002435952301324096,0.2009998778227811,0.3051792191603155,2.641649115599336),
    @(1.505080075141962,0.3397349683509674,0.115394263016715,0.6400460337125793,0.002437088981664155,0.1995071200968823,0.2972911854474489,2.655450574779024),
    @(1.475990901274315,0.3361132684547385,0.1157266018747166,0.6244449056556647,0.002437751803857276,0.1986575286035048,0.2927585722487152,2.66362851956417),
    @(1.466142038449675,0.3348866841923837,0.11584001966077,0.619163680173358,0.002437977777921999,0.1983714192014219,0.2912246659455562,2.666438559649023),
    @(1.510463926497664,0.3404051006751843,0.1153331789183798,0.642933953830422,0.002436967045595218,0.1996650952241623,0.2981304290009703,2.653956566754431),
    @(1.659388597920895,0.3589214352686554,0.113692376849407,0.7228739723491628,0.002433678913114657,0.2041191386024792,0.3213845898687353,2.614885904627215),
    @(1.756718253755992,0.3710030387857444,0.1126662570500458,0.7751780083420101,0.002431609932066052,0.2071099834488592,0.3366198862710732,2.591494959318766),
    @(1.704772213044635,0.3645567700389165,0.1132096891103909,0.7472568307830727,0.002432706887040349,0.205506419755082,0.3284851806874372,2.603782517392204),
    @(1.50802992485427,0.3401021452320663,0.1153607783775119,0.6416283278902171,0.002437022143814601,0.199593647832117,0.2977509997391508,2.654631250005849),
    @(1.295858068186362,0.313647599553434,0.1178799601580022,0.5279251897347166,0.002440291673589723,0.195546503624783,0.2647673283617635,2.71863887354904)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $colLetter = $cols[$j]
        $ws.Range("$colLetter$rowNum").Value = $rowValues[$j]
    }
}
